$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the existing "Arrow Order" entry (row 10) to "Arrow Order 1",
#    since a second Arrow order is being added below. The existing
#    hyperlink (rId9) is left pointing at its original target file.
$ws.Range("D10").Value = "Arrow Order 1.pdf"
$ws.Range("A10").Value = "Arrow Order 1"

# 2. Add the new "Arrow Order 2" receipt row (row 15), just above the Total row.
$ws.Range("A15").Value = "Arrow Order 2"

$ws.Range("B15").Value = 43117
$ws.Range("B15").NumberFormat = "mm/dd/yy;@"

$ws.Range("C15").Value = "Brian"

$ws.Hyperlinks.Add($ws.Range("D15"), "Arrow Order 2.pdf") | Out-Null
$ws.Range("D15").Value = "Arrow Order 2.pdf"
$ws.Range("D15").Style = "Hyperlink"

$ws.Range("E15").Value = 24.74
$ws.Range("E15").NumberFormat = """$""#,##0.00"

$ws.Range("F15").Value = "ATtiny85, Voltage Regulators, Rotary Position Sensors, MOSFETs, SMD Capacitors"

# 3. Refresh the total so it recalculates with the new row included.
$ws.Range("E16").Formula = "=SUM(E2:E15)"

$ws.Range("F16").Select()
